{"js": "// The document contains a single table whose 5 \"data\" rows (interleaved\n// with blank spacer rows) hold three-digit-by-one-digit multiplication\n// problems and their answers, e.g. \"202\u00d76=1212\". This edit swaps each\n// problem/answer pair for a new one, keeping cell formatting untouched.\n//\n// Target rows/cells are addressed by table position (row index, column\n// index) rather than by searching for the old text, because one value\n// (\"530\u00d79=4770\") occurs twice in the table and a text-based replace could\n// not tell the two occurrences apart.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// 0-based table row index -> the 5 new cell values (left to right) that\n// replace that row's current problem/answer strings.\nconst rowUpdates = [\n  { row: 0, values: [\"590\u00d79=5310\", \"271\u00d72=542\", \"350\u00d73=1050\", \"842\u00d74=3368\", \"603\u00d75=3015\"] },\n  { row: 4, values: [\"163\u00d72=326\", \"251\u00d76=1506\", \"204\u00d79=1836\", \"526\u00d73=1578\", \"178\u00d77=1246\"] },\n  { row: 9, values: [\"118\u00d79=1062\", \"923\u00d75=4615\", \"608\u00d77=4256\", \"156\u00d75=780\", \"131\u00d77=917\"] },\n  { row: 14, values: [\"882\u00d78=7056\", \"186\u00d76=1116\", \"155\u00d75=775\", \"795\u00d74=3180\", \"571\u00d76=3426\"] },\n  { row: 19, values: [\"314\u00d72=628\", \"546\u00d73=1638\", \"979\u00d72=1958\", \"749\u00d73=2247\", \"760\u00d78=6080\"] },\n];\n\nfor (const { row, values } of rowUpdates) {\n  for (let col = 0; col < values.length; col++) {\n    table.getCell(row, col).value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table whose 5 \"data\" rows (interleaved\n# with blank spacer rows) hold three-digit-by-one-digit multiplication\n# problems and their answers, e.g. \"202x6=1212\". This edit swaps each\n# problem/answer pair for a new one, keeping cell formatting untouched.\n#\n# Target rows/cells are addressed by table position (row, column) rather\n# than by searching/replacing the old text, because one value\n# (\"530x9=4770\") occurs twice in the table and a text-based replace could\n# not tell the two occurrences apart.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowUpdates = @{\n    1  = @(\"590\u00d79=5310\", \"271\u00d72=542\", \"350\u00d73=1050\", \"842\u00d74=3368\", \"603\u00d75=3015\")\n    5  = @(\"163\u00d72=326\", \"251\u00d76=1506\", \"204\u00d79=1836\", \"526\u00d73=1578\", \"178\u00d77=1246\")\n    10 = @(\"118\u00d79=1062\", \"923\u00d75=4615\", \"608\u00d77=4256\", \"156\u00d75=780\", \"131\u00d77=917\")\n    15 = @(\"882\u00d78=7056\", \"186\u00d76=1116\", \"155\u00d75=775\", \"795\u00d74=3180\", \"571\u00d76=3426\")\n    20 = @(\"314\u00d72=628\", \"546\u00d73=1638\", \"979\u00d72=1958\", \"749\u00d73=2247\", \"760\u00d78=6080\")\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $values = $rowUpdates[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
